$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.659.07"
$ws.Range("E2").Value = "  +7.68%  "

$ws.Range("D3").Value = "3.590.60"
$ws.Range("E3").Value = "  +3.46%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.02"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.43"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("E7").Value = "  +3.63%  "

$ws.Range("D8").Value = "3.581.92"
$ws.Range("E8").Value = "  +3.40%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.773"
$ws.Range("E10").Value = "  +6.35%  "

$ws.Range("E11").Value = "  +17.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000336"
$ws.Range("E12").Value = "  +44.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.32"
$ws.Range("E13").Value = "  -0.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.90"
$ws.Range("E14").Value = "  +2.17%  "

$ws.Range("D15").Value = "4.153.80"
$ws.Range("E15").Value = "  +3.28%  "

$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.40"
$ws.Range("E17").Value = "  -0.83%  "

$ws.Range("D18").Value = "3.605.78"
$ws.Range("E18").Value = "  +3.28%  "

$ws.Range("E19").Value = "  +5.85%  "

$ws.Range("D20").Value = "67.487.35"
$ws.Range("E20").Value = "  +7.57%  "

$ws.Range("E21").Value = "  -2.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "452.36"
$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.08"
$ws.Range("E23").Value = "  -1.54%  "

$ws.Range("E24").Value = "  -4.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.14"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").Value = "  -6.95%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.00"
$ws.Range("E28").Value = "  +4.73%  "

$ws.Range("E29").Value = "  +2.04%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.78"
$ws.Range("E30").Value = "  +3.86%  "

$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.35"
$ws.Range("E31").Value = "  +2.13%  "

$ws.Range("E32").Value = "  +4.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.40"
$ws.Range("E33").Value = "  -2.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.162"
$ws.Range("E34").Value = "  -3.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.75"
$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.80"
$ws.Range("E37").Value = "  -2.68%  "

$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("D39").Value = "0.0₃0770"
$ws.Range("E39").Value = "  +36.75%  "

$ws.Range("E40").Value = "  +9.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("E41").Value = "  -0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.06"
$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "149.35"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("E45").Value = "  -2.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.26"
$ws.Range("E46").Value = "  -2.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.31"
$ws.Range("E47").Value = "  -1.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.98"
$ws.Range("E48").Value = "  -4.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -2.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "115.79"
$ws.Range("E50").Value = "  +6.13%  "

$ws.Range("E51").Value = "  +10.96%  "
